$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value in E10 (num_matches changed for season index 8)
$ws.Range("E10").Value = 854877

# Add new row 11 for season M2_10 Cat 2020
# Copy formatting from the row above (row 10) for the new row, then set values
$ws.Range("A10:H10").Copy() | Out-Null
$ws.Range("A11:H11").PasteSpecial(-4122) | Out-Null

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "M2_10 Cat 2020"
$ws.Range("C11").Value = 9703
$ws.Range("D11").Value = 10804
$ws.Range("E11").Value = 929613
$ws.Range("F11").Value = 9977
$ws.Range("G11").Value = 10067
$ws.Range("H11").Value = 10176
